$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the header row text (confusion-matrix column headers renamed,
#    a new "Test numbers " header added in F1 above the existing SUM
#    formula column, and a new "Model Version" header added last in A1).
$ws.Range("F1").Value = "Test numbers "
$ws.Range("B1").Value = "Actual 0 _Predicted 0"
$ws.Range("C1").Value = "Actual 1_Predicted 1"
$ws.Range("D1").Value = "Actual 0_Predicted 1"
$ws.Range("E1").Value = "Actual 1_Predicted 0"
$ws.Range("G1").Value = "Error"
$ws.Range("A1").Value = "Model Version"

# 2. Remove the old scratch columns (I: extra sum, J: extra ratio, K: extra
#    text labels) that are no longer part of the confusion matrix.
$ws.Range("I1:K10").ClearContents()

# 3. Turn the matrix into a proper Excel Table ("Table5") with a table style.
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:G10"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table5"
$lo.TableStyle = "TableStyleMedium20"

# 4. Give the "Error" column its percentage display style.
$ws.Range("G2:G10").Style = "Percent"

$wb.Save()
